$d = $word.ActiveDocument

# The edit appends 8 blank paragraphs followed by one paragraph containing
# new text to the last cell of the last row of the second table (the cell
# ending in "...then tunes and adjusts its behavior accordingly.").
$t = $d.Tables.Item(2)
$cell = $t.Cell($t.Rows.Count, 1)
$r = $cell.Range
$r.Collapse(0)

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$emptyPara = '<w:p xmlns:w="' + $w + '"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'

$textPara = '<w:p xmlns:w="' + $w + '"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Lets read this again n again about agile just checking github delte</w:t></w:r></w:p>'

$combined = ""
for ($i = 0; $i -lt 8; $i++) {
    $combined = $combined + $emptyPara
}
$combined = $combined + $textPara

$r.InsertXML($combined)
